$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of the row-pair (rows 2-3) with the row-pair
# (rows 4-5): row2<->row4 content and row3<->row5 content (the two
# records that used to be first in the sheet move down, and the two
# that used to be last move up), while row 1 (header) is untouched.
#
# Range.Copy is used (instead of assigning literal strings through
# .Value) so each cell's original literal representation is preserved
# verbatim - text that looks like a date (e.g. "2015-07-02") is not
# reinterpreted/reformatted as a date serial number the way a plain
# string assignment would.
#
# Because Range.Copy silently skips writing into destination cells
# whose corresponding source cell is blank (it does not clear them),
# every destination is explicitly cleared with ClearContents() right
# before the paste that is meant to populate it, so columns that are
# populated in one row-pair but empty in the other end up correctly
# empty/populated after the swap.

$topRange = $ws.Range("A2:AY3")
$bottomRange = $ws.Range("A4:AY5")
$scratchRange = $ws.Range("A1000:AY1001")

# 1) stash the bottom pair (rows 4-5) in a scratch area far away
$scratchRange.ClearContents()
$bottomRange.Copy($scratchRange)

# 2) move the top pair (rows 2-3) down into rows 4-5
$bottomRange.ClearContents()
$topRange.Copy($bottomRange)

# 3) move the stashed original bottom pair up into rows 2-3
$topRange.ClearContents()
$scratchRange.Copy($topRange)

# 4) clean up the scratch area
$scratchRange.ClearContents()
